# Update attendance/visitor figures (column F) on the "展览" and "全部类型"
# sheets to reflect the latest generated output.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1468
$ws1.Range("F16").Value = 70
$ws1.Range("F18").Value = 4740
$ws1.Range("F20").Value = 823
$ws1.Range("F22").Value = 2202

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1468
$ws4.Range("F16").Value = 70
$ws4.Range("F18").Value = 4740
$ws4.Range("F22").Value = 823
$ws4.Range("F24").Value = 2202
